$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row's updated Price (D) and Volume(1h) (E) values.
# These columns hold plain text in the source workbook (e.g. "301.79",
# "-1.19%"), so we force text entry to avoid Excel auto-converting them
# into numeric/percentage cells.
$updates = @(
    @{ Row = 2;  D = "300.70";     E = "-1.82%" }
    @{ Row = 3;  D = "31.53";      E = "-2.44%" }
    @{ Row = 4;  D = "5.152";      E = "-2.82%" }
    @{ Row = 5;  D = "0.07493";    E = "0.91%" }
    @{ Row = 6;  D = "7.842";      E = "0.41%" }
    @{ Row = 7;  D = "3.826";      E = "3.43%" }
    @{ Row = 8;  D = "1.656";      E = "1.64%" }
    @{ Row = 9;  D = "0.9241";     E = "0.65%" }
    @{ Row = 10; D = "0.1712";     E = "2.33%" }
    @{ Row = 11; D = "0.07695";    E = "5.05%" }
    @{ Row = 12; D = "0.08023";    E = "-0.43%" }
    @{ Row = 13; D = "0.02998";    E = "-3.65%" }
    @{ Row = 14; D = "0.09899";    E = "0.53%" }
    @{ Row = 15; D = "0.001490";   E = "-1.89%" }
    @{ Row = 16; D = "0.04667";    E = "2.77%" }
    @{ Row = 17; D = "0.006187";   E = "1.71%" }
    @{ Row = 18; D = "3.445";      E = "-1.22%" }
    @{ Row = 19; D = "2.231";      E = "-0.75%" }
    @{ Row = 20; D = "0.3293";     E = "0.62%" }
    @{ Row = 21; D = "0.1337";     E = "2.16%" }
    @{ Row = 22; D = "4.584";      E = "7.70%" }
    @{ Row = 23; D = "0.1552";     E = "-5.14%" }
    @{ Row = 24; D = "0.001223";   E = "-0.01%" }
    @{ Row = 25; D = "0.004425";   E = "-2.46%" }
    @{ Row = 26; D = $null;        E = "20.03%" }
    @{ Row = 27; D = "0.0001799";  E = "5.33%" }
    @{ Row = 39; D = "0.01661";    E = "0.73%" }
    @{ Row = 40; D = "0.04542";    E = "0.34%" }
    @{ Row = 41; D = "0.006971";   E = "-4.31%" }
    @{ Row = 42; D = "0.1343";     E = "-1.78%" }
    @{ Row = 43; D = $null;        E = "-5.21%" }
    @{ Row = 44; D = "0.01236";    E = "-12.48%" }
    @{ Row = 45; D = "0.00006044"; E = "1.65%" }
    @{ Row = 46; D = "1.930";      E = "1.95%" }
    @{ Row = 47; D = "0.01226";    E = "-5.47%" }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cellD = $ws.Range("D$r")
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }
    $cellE = $ws.Range("E$r")
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
}
